$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46025
$ws.Range("B2").Value = 80.37
$ws.Range("C2").Value = 78.69
$ws.Range("D2").Value = 76.86
$ws.Range("E2").Value = 75.88
$ws.Range("F2").Value = 73.94
$ws.Range("G2").Value = 75.27
$ws.Range("H2").Value = 76.56999999999999
$ws.Range("I2").Value = 79.77
$ws.Range("J2").Value = 90.69
$ws.Range("K2").Value = 94.09
$ws.Range("L2").Value = 93.31
$ws.Range("M2").Value = 88.70999999999999
$ws.Range("N2").Value = 87.41
$ws.Range("O2").Value = 86.61
$ws.Range("P2").Value = 89.83
$ws.Range("Q2").Value = 95.09
$ws.Range("R2").Value = 100.21
$ws.Range("S2").Value = 108.09
$ws.Range("T2").Value = 103.78
$ws.Range("U2").Value = 101.02
$ws.Range("V2").Value = 98.03
$ws.Range("W2").Value = 91.73
$ws.Range("X2").Value = 90.40000000000001
$ws.Range("Y2").Value = 85.95
$ws.Range("Z2").Value = 88.43000000000001
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 103.28
$ws.Range("AC2").Value = "16h-18h"
$ws.Range("AD2").Value = 104.15
$ws.Range("AF2").Value = 102.4
$ws.Range("AG2").Value = "0h-23h"
